$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1) changed
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2:E2) changed
$ws.Range("B2").Value = 19.700776789929254
$ws.Range("C2").Value = 14.336164359257241
$ws.Range("D2").Value = 27.207993566727964
$ws.Range("E2").Value = 14.999933945741793

# Row 3 data values (B3:E3) changed
$ws.Range("B3").Value = 26.675053599404535
$ws.Range("C3").Value = 15.606197544409987
$ws.Range("D3").Value = 30.344649674770377
$ws.Range("E3").Value = 13.824821612113388

# Selection narrowed from B1:AY3 to B1:E3
$ws.Range("B1:E3").Select()
